$d = $word.ActiveDocument

# Paragraph 1: "Assunto: " -> "Assunto: Retificação de área<TAB><TAB>1ª análise"
$d.Paragraphs(1).Range.Text = "Assunto: Retificação de área`t`t1ª análise"

# Paragraph 2: "Solicitação de demanda: " -> "...Processo Físico"
$d.Paragraphs(2).Range.Text = "Solicitação de demanda: Processo Físico"

# Paragraph 3: "Contribuinte: " -> "Contribuinte: Carlos"
$d.Paragraphs(3).Range.Text = "Contribuinte: Carlos"

# Paragraph 4 was "Endereço do imóvel: " but the edited document swaps the
# order so that "Inscrição Imobiliária" now comes before "Endereço do imóvel".
$d.Paragraphs(4).Range.Text = "Inscrição Imobiliária: 123456789123456"

# Paragraph 5 was "Inscrição Imobiliária: " -> now holds the address text.
$d.Paragraphs(5).Range.Text = "Endereço do imóvel: Rua Cibele, nº 145 - bairro Vila Amelioa, Itabira - MG"

# Paragraph 6: "Dados recebidos: " -> adds a list item plus an embedded
# newline character right before the end of the run.
$d.Paragraphs(6).Range.Text = "Dados recebidos: - Planta do imóvel`n"

# Add a new, empty trailing paragraph (matching the same run formatting)
$d.Paragraphs(6).Range.InsertParagraphAfter()
$d.Paragraphs(7).Range.Text = ""
